$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (ALC)
$ws.Range("H15").Value = 1176.7347
$ws.Range("I15").Value = 1176.7347
$ws.Range("K15").Value = 3530.2041
$ws.Range("M15").Value = -3361.2041

# Row 41 (ALC)
$ws.Range("H41").Value = 791.2222
$ws.Range("J41").Value = 650
$ws.Range("L41").Value = 650
$ws.Range("N41").Value = -1530

# Row 98 (ALC)
$ws.Range("H98").Value = 1000.5
$ws.Range("I98").Value = 1077.0769
$ws.Range("J98").Value = 668.6667
$ws.Range("K98").Value = 1077.0769
$ws.Range("L98").Value = 668.6667
$ws.Range("M98").Value = 420.9231
$ws.Range("N98").Value = -3664.6667

# Row 99 (ALC)
$ws.Range("H99").Value = 716
$ws.Range("I99").Value = 688
$ws.Range("J99").Value = 800
$ws.Range("K99").Value = 2064
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -566
$ws.Range("N99").Value = -5396

# Row 122 (ALC)
$ws.Range("H122").Value = 1000.5
$ws.Range("I122").Value = 1077.0769
$ws.Range("J122").Value = 668.6667
$ws.Range("K122").Value = 3231.2307
$ws.Range("L122").Value = 2006.0001
$ws.Range("M122").Value = -781.2307000000001
$ws.Range("N122").Value = -6906.0001

# Row 123 (ALC)
$ws.Range("H123").Value = 43749.75
$ws.Range("J123").Value = 43749.75
$ws.Range("L123").Value = 43749.75
$ws.Range("N123").Value = -53549.75

# Row 124 (ALC)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()

# Row 125 (ALC)
$ws.Range("H125").Value = 15127.429
$ws.Range("I125").Value = 695
$ws.Range("J125").Value = 20900.4
$ws.Range("K125").Value = 6255
$ws.Range("L125").Value = 188103.6
$ws.Range("M125").Value = -3795
$ws.Range("N125").Value = -193023.6

# Row 127 (ALC)
$ws.Range("H127").Value = 38462470
$ws.Range("I127").Value = 579.4
$ws.Range("J127").Value = 47620064
$ws.Range("K127").Value = 1738.2
$ws.Range("L127").Value = 142860192
$ws.Range("M127").Value = 3221.8
$ws.Range("N127").Value = -142870112

# Row 131 (ALC)
$ws.Range("H131").Value = 5464.05
$ws.Range("I131").Value = 1356.5714
$ws.Range("J131").Value = 7675.769
$ws.Range("K131").Value = 4069.7142
$ws.Range("L131").Value = 23027.307
$ws.Range("M131").Value = 970.2857999999997
$ws.Range("N131").Value = -33107.307

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 5303.6377
$ws.Range("I32").Value = 4106.73
$ws.Range("K32").Value = 4106.73
$ws.Range("M32").Value = -3819.73

# Row 132 (ARM)
$ws.Range("H132").Value = 9413.414000000001
$ws.Range("I132").Value = 8653.929
$ws.Range("J132").Value = 10122.267
$ws.Range("K132").Value = 25961.787
$ws.Range("L132").Value = 30366.801
$ws.Range("M132").Value = -23431.787
$ws.Range("N132").Value = -35426.801

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 1684.5714
$ws.Range("I94").Value = 1548
$ws.Range("J94").Value = 1866.6666
$ws.Range("K94").Value = 1548
$ws.Range("L94").Value = 1866.6666
$ws.Range("M94").Value = -1097
$ws.Range("N94").Value = -2768.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 111.666664
$ws.Range("I22").Value = 84.166664
$ws.Range("J22").Value = 166.66667
$ws.Range("K22").Value = 84.166664
$ws.Range("L22").Value = 166.66667
$ws.Range("M22").Value = 265.833336
$ws.Range("N22").Value = -866.6666700000001

# Row 58 (CRP)
$ws.Range("H58").Value = 2220209
$ws.Range("I58").Value = 4330839.5
$ws.Range("J58").Value = 4047.1
$ws.Range("K58").Value = 4330839.5
$ws.Range("L58").Value = 4047.1
$ws.Range("M58").Value = -4330636.5
$ws.Range("N58").Value = -4453.1

# Row 136 (CRP)
$ws.Range("H136").Value = 2220209
$ws.Range("I136").Value = 4330839.5
$ws.Range("J136").Value = 4047.1
$ws.Range("K136").Value = 12992518.5
$ws.Range("L136").Value = 12141.3
$ws.Range("M136").Value = -12989968.5
$ws.Range("N136").Value = -17241.3

$ws = $wb.Worksheets.Item("CUL")
# Row 22 (CUL)
$ws.Range("H22").Value = 100001640
$ws.Range("J22").Value = 2177.4285
$ws.Range("L22").Value = 6532.2855
$ws.Range("N22").Value = -6870.2855

# Row 27 (CUL)
$ws.Range("H27").Value = 100001640
$ws.Range("J27").Value = 2177.4285
$ws.Range("L27").Value = 6532.2855
$ws.Range("N27").Value = -6736.2855

# Row 36 (CUL)
$ws.Range("H36").Value = 975
$ws.Range("I36").Value = 633.3333
$ws.Range("K36").Value = 1899.9999
$ws.Range("M36").Value = -1730.9999

# Row 54 (CUL)
$ws.Range("H54").Value = 2486.5
$ws.Range("I54").Value = 2000
$ws.Range("J54").Value = 2583.8
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 7751.400000000001
$ws.Range("M54").Value = -5441
$ws.Range("N54").Value = -8869.400000000001

# Row 131 (CUL)
$ws.Range("H131").Value = 653
$ws.Range("I131").Value = 327.66666
$ws.Range("J131").Value = 710.41174
$ws.Range("K131").Value = 982.9999799999999
$ws.Range("L131").Value = 2131.23522
$ws.Range("M131").Value = 4057.00002
$ws.Range("N131").Value = -12211.23522

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 1230.9166
$ws.Range("I97").Value = 1007.7857
$ws.Range("J97").Value = 1543.3
$ws.Range("K97").Value = 1007.7857
$ws.Range("L97").Value = 1543.3
$ws.Range("M97").Value = -511.7857
$ws.Range("N97").Value = -2535.3

# Row 123 (GSM)
$ws.Range("H123").Value = 16540
$ws.Range("J123").Value = 16540
$ws.Range("L123").Value = 16540
$ws.Range("N123").Value = -21440

# Row 124 (GSM)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 1925.174
$ws.Range("I132").Value = 1397.25
$ws.Range("J132").Value = 2501.0908
$ws.Range("K132").Value = 4191.75
$ws.Range("L132").Value = 7503.2724
$ws.Range("M132").Value = -1661.75
$ws.Range("N132").Value = -12563.2724

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 2340.1428
$ws.Range("I16").Value = 2400.1667
$ws.Range("J16").Value = 1980
$ws.Range("K16").Value = 2400.1667
$ws.Range("L16").Value = 1980
$ws.Range("M16").Value = -2230.1667
$ws.Range("N16").Value = -2320

# Row 22 (LTW)
$ws.Range("H22").Value = 336.75
$ws.Range("I22").Value = 304.1
$ws.Range("K22").Value = 304.1
$ws.Range("M22").Value = -9.100000000000023

# Row 27 (LTW)
$ws.Range("H27").Value = 336.75
$ws.Range("I27").Value = 304.1
$ws.Range("K27").Value = 304.1
$ws.Range("M27").Value = -197.1

# Row 32 (LTW)
$ws.Range("H32").Value = 520.4286
$ws.Range("I32").Value = 520.4286
$ws.Range("K32").Value = 520.4286
$ws.Range("M32").Value = -203.4286

# Row 55 (LTW)
$ws.Range("H55").Value = 200209.25
$ws.Range("I55").Value = 400160.1
$ws.Range("J55").Value = 258.4
$ws.Range("K55").Value = 400160.1
$ws.Range("L55").Value = 258.4
$ws.Range("M55").Value = -399987.1
$ws.Range("N55").Value = -604.4
